# BulkTool: the "Datum" field is now a free-form string instead of a
# strict DateOnly value, so replace the ISO date that was stored for
# François Dupont (row 6, column C = "Datum") with a looser text value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value = "December 2024"

# Column A (the "Naam" column) becomes very slightly narrower, and
# column C (the "Datum" column, now holding longer free-text values
# like "December 2024") gets its own explicit width.
$ws.Range("A:A").ColumnWidth = 25.75
$ws.Range("C:C").ColumnWidth = 12.6

# Move the selection / active cell.
$ws.Range("C7").Select()
